# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet.
$wb = $excel.ActiveWorkbook
$snapshot = $wb.Worksheets.Item("snapshot")

$updates = @(
    @{ Row = 2; Value = "2025-12-02T07:01:46.757060+00:00" },
    @{ Row = 3; Value = "2025-12-02T07:01:49.495128+00:00" },
    @{ Row = 4; Value = "2025-12-02T07:01:49.495157+00:00" },
    @{ Row = 5; Value = "2025-12-02T07:01:51.807745+00:00" },
    @{ Row = 6; Value = "2025-12-02T07:01:54.174149+00:00" },
    @{ Row = 7; Value = "2025-12-02T07:01:57.095415+00:00" },
    @{ Row = 8; Value = "2025-12-02T07:01:59.400570+00:00" },
    @{ Row = 9; Value = "2025-12-02T07:02:04.068314+00:00" },
    @{ Row = 10; Value = "2025-12-02T07:02:04.068342+00:00" },
    @{ Row = 11; Value = "2025-12-02T07:02:06.440100+00:00" },
    @{ Row = 12; Value = "2025-12-02T07:02:08.811888+00:00" },
    @{ Row = 13; Value = "2025-12-02T07:02:08.811919+00:00" },
    @{ Row = 14; Value = "2025-12-02T07:02:11.097661+00:00" },
    @{ Row = 15; Value = "2025-12-02T07:02:13.840177+00:00" },
    @{ Row = 16; Value = "2025-12-02T07:02:13.840205+00:00" },
    @{ Row = 17; Value = "2025-12-02T07:02:16.632313+00:00" },
    @{ Row = 18; Value = "2025-12-02T07:02:19.436433+00:00" },
    @{ Row = 19; Value = "2025-12-02T07:02:19.436463+00:00" },
    @{ Row = 20; Value = "2025-12-02T07:02:22.182014+00:00" },
    @{ Row = 21; Value = "2025-12-02T07:02:22.182045+00:00" },
    @{ Row = 22; Value = "2025-12-02T07:02:22.182063+00:00" },
    @{ Row = 23; Value = "2025-12-02T07:02:24.540313+00:00" },
    @{ Row = 24; Value = "2025-12-02T07:02:29.685947+00:00" },
    @{ Row = 25; Value = "2025-12-02T07:02:32.071080+00:00" },
    @{ Row = 26; Value = "2025-12-02T07:02:32.071113+00:00" },
    @{ Row = 27; Value = "2025-12-02T07:02:34.855661+00:00" },
    @{ Row = 28; Value = "2025-12-02T07:02:34.855691+00:00" },
    @{ Row = 29; Value = "2025-12-02T07:02:37.667374+00:00" },
    @{ Row = 30; Value = "2025-12-02T07:02:37.667402+00:00" }
)

foreach ($u in $updates) {
    $snapshot.Cells.Item($u.Row, 11).Value = $u.Value
}

# The "new_injured" sheet's two pending rows (Torpedo / Traktor) were
# processed, so remove them from the worksheet, shrinking it back to
# just the header row.
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows("2:3").Delete()

$wb.Save()
